# Swap the "title" and "uri" values between row 2 and row 3 (Japanese
# typhoon article row <-> WPC Surface Analysis Archive row), leaving the
# timestamp / historical distance / time bucket columns untouched since
# they already hold identical values in both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleRow2 = $ws.Range("A2").Value()
$titleRow3 = $ws.Range("A3").Value()
$uriRow2   = $ws.Range("E2").Value()
$uriRow3   = $ws.Range("E3").Value()

$ws.Range("A2").Value = $titleRow3
$ws.Range("A3").Value = $titleRow2

$ws.Range("E2").Value = $uriRow3
$ws.Range("E3").Value = $uriRow2
